# Update OpenCart Database Test Cases.xlsx
# - Rename test case IDs from "DBQ-NN" to "TC_DBQ_0NN" (rows 2-24, column A)
# - Widen column A to fit the new, longer ID strings

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ids = @(
    "TC_DBQ_001",
    "TC_DBQ_002",
    "TC_DBQ_003",
    "TC_DBQ_004",
    "TC_DBQ_005",
    "TC_DBQ_006",
    "TC_DBQ_007",
    "TC_DBQ_008",
    "TC_DBQ_009",
    "TC_DBQ_010",
    "TC_DBQ_011",
    "TC_DBQ_012",
    "TC_DBQ_013",
    "TC_DBQ_014",
    "TC_DBQ_015",
    "TC_DBQ_016",
    "TC_DBQ_017",
    "TC_DBQ_018",
    "TC_DBQ_019",
    "TC_DBQ_020",
    "TC_DBQ_021",
    "TC_DBQ_022",
    "TC_DBQ_023"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}

# Widen column A (was 9.57 chars) to fit the longer "TC_DBQ_0NN" ids (now 15.71 chars)
$ws.Columns.Item(1).ColumnWidth = 14.87
